# Applies the "flujos alterno del modulo de autenticacion" edit to the
# autenticacion.xlsx workbook (sheet "Datos").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")
$ws.Activate()

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Column G: rename the test-case codes.
#   G2: SVPPRU08  -> OSVPPRU10
#   G3: (empty)   -> OSVPPRU10
#   G6: USUARIOS41-> OSVPPRU04
# ---------------------------------------------------------------------
$ws.Range("G2").Value = "OSVPPRU10"
$ws.Range("G3").Value = "OSVPPRU10"
$ws.Range("G6").Value = "OSVPPRU04"

# ---------------------------------------------------------------------
# Column F: add the missing document number for row 8.
#   F8: (empty) -> 854124014
# ---------------------------------------------------------------------
$ws.Range("F8").Value = "854124014"

# ---------------------------------------------------------------------
# Column E: append a trailing period to the "usuario o clave invalida"
# error message (rows 5, 6 and 8).
# ---------------------------------------------------------------------
$ws.Range("E5").Value = "Usuario o clave inválida. Inténtalo nuevamente."
$ws.Range("D5").Copy()
$ws.Range("E5").PasteSpecial($xlPasteFormats)

$ws.Range("E6").Value = "Usuario o clave inválida. Inténtalo nuevamente."
$ws.Range("D6").Copy()
$ws.Range("E6").PasteSpecial($xlPasteFormats)

$ws.Range("E8").Value = "Usuario o clave inválida. Inténtalo nuevamente."
$ws.Range("D8").Copy()
$ws.Range("E8").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Column E: add the new greeting message for the new alternate-flow
# rows (2 and 3).
# ---------------------------------------------------------------------
$ws.Range("E2").Value = "¡Hola!"
$ws.Range("C2").Copy()
$ws.Range("E2").PasteSpecial($xlPasteFormats)

$ws.Range("E3").Value = "¡Hola!"
$ws.Range("C3").Copy()
$ws.Range("E3").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Complete row 3 (new alternate test case) and clear the obsolete
# test-case code in row 8.
# ---------------------------------------------------------------------
$ws.Range("H3").Value = 1234
$ws.Range("G8").ClearContents()

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Sheet view: drop the B1 frozen/top-left scroll position and move the
# active selection from E12 to E3.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("E3").Select()
